$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.021.41'
$ws.Range("E2").Value = '  +4.53%  '

$ws.Range("D3").Value = '2.286.72'
$ws.Range("E3").Value = '  +4.98%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.96'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.42%  '

$ws.Range("E6").Value = '  +3.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.91'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +9.73%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  +12.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.64'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +9.12%  '

$ws.Range("E11").Value = '  +4.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.81'
$ws.Range("D12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.50'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +9.23%  '

$ws.Range("E14").Value = '  +1.95%  '

$ws.Range("D15").Value = '2.628.87'
$ws.Range("E15").Value = '  +4.97%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.14'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +5.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.897'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +6.08%  '

$ws.Range("D18").Value = '2.301.00'
$ws.Range("E18").Value = '  +6.02%  '

$ws.Range("D19").Value = '42.943.98'
$ws.Range("E19").Value = '  +4.49%  '

$ws.Range("E20").Value = '  +6.68%  '

$ws.Range("E21").Value = '  +6.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.69'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.88%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.05'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.45%  '

$ws.Range("E24").Value = '  +8.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.91'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.85'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.87%  '

$ws.Range("E27").Value = '  -0.19%  '

$ws.Range("E28").Value = '  +1.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.69'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.14'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.22'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +5.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.35'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +12.44%  '

$ws.Range("E34").Value = '  +6.79%  '

$ws.Range("E35").Value = '  +8.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.42'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +27.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.85'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +21.76%  '

$ws.Range("E38").Value = '  +5.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.79'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.11%  '

$ws.Range("E40").Value = '  +2.53%  '

$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.47'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +19.54%  '

$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.34'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.12'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +10.50%  '

$ws.Range("E44").Value = '  +13.09%  '

$ws.Range("E45").Value = '  +8.74%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.95'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -9.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '61.75'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.79%  '

$ws.Range("E48").Value = '  +4.74%  '

$ws.Range("E49").Value = '  +5.04%  '

$ws.Range("E50").Value = '  +0.22%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.87'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.21%  '
